$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update from the diff. Every new value is written with a
# leading apostrophe (quote-prefix) so Excel always stores it as literal text
# -- this matters for numeric-looking strings such as "53.822.44" or "0.0200"
# which must stay text (not be coerced to a number) to match the source data.
# The style is then reset to Normal so the quote-prefix / text-format marker
# doesn't leave a stray style on the cell (keeping formatting identical to
# untouched neighbouring cells).

$ws.Range('D2').Value = '''53.822.44'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -4.56%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''2.226.46'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -6.38%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  -0.04%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''484.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -4.26%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''124.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -4.46%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.23%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = '''0.521'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -4.45%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''2.225.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -6.80%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''0.0920'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -7.04%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.149'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -1.01%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = '''Cardano'
$ws.Range('B12').Style = 'Normal'
$ws.Range('C12').Value = '''https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('C12').Style = 'Normal'
$ws.Range('D12').Value = '''0.315'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -3.59%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = '''Toncoin'
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = '''4.61'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -4.57%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''2.616.02'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -6.62%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''21.08'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -2.05%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''53.495.64'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -5.05%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''0.0000127'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -4.22%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''2.223.42'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -6.34%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value = '''Chainlink'
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').Value = '''9.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -5.00%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('B20').Value = '''Polkadot'
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').Value = '''3.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -2.60%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''293.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -4.95%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = '''6.13'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -2.45%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -0.04%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''62.78'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -5.23%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = '''  +0.23%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''0.367'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -0.13%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '''  -0.94%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''2.304.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -7.32%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''7.02'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -3.43%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''164.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -5.02%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''1.57'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -4.71%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = '''USDe'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = '''https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = '''0.999'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -0.03%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = '''PEPE'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = '''https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = '''0.0₃0668'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -6.34%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''5.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -1.33%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  +0.30%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''1.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -2.70%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''17.28'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -2.48%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''1.16'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -1.70%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.842'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  +2.42%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''3.56'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -5.15%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''35.07'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -3.98%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.368'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +0.26%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''1.36'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -1.65%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''3.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -2.86%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = '''Aave'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = '''125.03'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -3.20%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = '''RenderToken'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = '''https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''4.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -4.08%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''0.0877'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -2.36%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''0.532'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -6.60%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = '''233.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -2.99%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''0.0470'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -2.54%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''0.0200'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -3.65%  '
$ws.Range('E51').Style = 'Normal'
